$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (shifts existing rows 3..53 down to 4..54,
# carrying formatting down with them, matching the target diff).
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with the new weekly price-record values.
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = "Vega Monumental Concepción"
$ws.Range("C3").Value = "Bíobío"
$ws.Range("D3").Value = 45083
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 100114007
$ws.Range("G3").Value = "Jengibre"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 30
$ws.Range("K3").Value = 18000
$ws.Range("L3").Value = 18000
$ws.Range("M3").Value = 18000
$ws.Range("N3").Value = "$/caja 13 kilos"
$ws.Range("O3").Value = "Perú"
$ws.Range("P3").Value = 1385
$ws.Range("Q3").Value = 13
$ws.Range("R3").Value = "Hortaliza"
